$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A14:H14").ClearContents()
$ws.Rows.Item(14).AutoFit()
